$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    "heilen",
    "drucken",
    "mühen",
    "ruhen",
    "scheitern",
    "grüßen",
    "stürzen",
    "stecken",
    "nerven",
    "schwächen",
    "hoffen",
    "wundern",
    "wüten",
    "warnen",
    "helfen",
    "stürmen",
    "wählen",
    "sorgen",
    "gründen",
    "kratzen",
    "platzen",
    "handeln",
    "gelten",
    "lesen",
    "sprengen",
    "trauen",
    "bauen",
    "wenden",
    "herzen",
    "kehren",
    "fließen",
    "hacken",
    "melden",
    "lehnen",
    "weigern",
    "werden",
    "husten",
    "führen",
    "albern",
    "äußern",
    "wehren",
    "pfeifen",
    "dringen",
    "ehren",
    "fehlen",
    "sinken",
    "freuen",
    "flehen",
    "zögern",
    "wachsen",
    "kümmern",
    "schnellen",
    "süßen",
    "arten",
    "wirken",
    "betteln",
    "heben",
    "krachen",
    "fahren",
    "klettern",
    "treiben",
    "schulden",
    "schicken",
    "schenken",
    "reizen",
    "werfen",
    "münzen",
    "drohen",
    "fallen",
    "filmen",
    "achten",
    "spannen",
    "hören",
    "kennen",
    "ärgern",
    "leihen",
    "trennen",
    "fordern",
    "flüchten",
    "rufen",
    "wiegen",
    "geben",
    "zeugen",
    "spinnen",
    "quälen",
    "tanzen",
    "graben",
    "dienen",
    "liegen",
    "tropfen",
    "tauschen",
    "beten",
    "knarren",
    "lassen",
    "schwingen",
    "seufzen",
    "flüstern",
    "spielen",
    "lockern",
    "pflanzen",
    "ändern",
    "zeichnen",
    "fällen",
    "mauern",
    "kichern",
    "boxen",
    "atmen",
    "runden",
    "scheinen",
    "zielen",
    "stoßen",
    "drehen",
    "enden",
    "fischen",
    "stärken",
    "kosten",
    "mögen",
    "sperren",
    "kriegen",
    "schrecken",
    "backen",
    "schätzen",
    "liefern",
    "fühlen",
    "stellen",
    "folgen",
    "greifen",
    "irren",
    "narren",
    "jubeln",
    "stammen",
    "erben",
    "schwören",
    "weichen",
    "hauen",
    "morden",
    "formen",
    "binden",
    "streichen",
    "biegen",
    "rasen",
    "brauchen",
    "dürfen",
    "schreiten",
    "trösten",
    "schneiden",
    "schauen",
    "schlucken",
    "feiern",
    "sammeln",
    "fesseln",
    "bluten",
    "meinen",
    "zeigen",
    "prügeln",
    "tollen",
    "zünden",
    "fangen",
    "räumen",
    "wehen",
    "heulen",
    "bitten",
    "suchen",
    "sichern",
    "knurren",
    "streifen",
    "altern",
    "schmecken",
    "öffnen",
    "leisten",
    "malen",
    "bergen",
    "saufen",
    "lieben",
    "sterben",
    "bellen",
    "stören",
    "wehtun",
    "kranken",
    "siegen",
    "stehlen",
    "jagen",
    "machen",
    "planen",
    "rauchen",
    "klingen",
    "stillen",
    "decken",
    "loben",
    "lügen",
    "töten",
    "spüren"
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $words[$i]
}

